# Veno-Echo BOM.xlsx corrections:
#  - Row 26: was a stray/incorrect "CONN IC DIP SOCKET 40POS GOLD / 40 pin
#    socket / ED90044-ND / 575-11043640" line; corrected to the actual
#    "20 Pin Female header socket(s)" part with qty 2 and no Digikey/Mouser
#    part numbers (not SMD pre-populated).
#  - Row 36: Thonkiconn jack part number/description corrected from the
#    stereo "OG-PJ301_PJ301_THONKICONN6" / "Thonkiconn" to the mono
#    variant " PJ301M-12 / PJ398SM" / "mono Thonkiconn".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: 20 pin female header sockets ---
$ws.Range("A26").Value = 2
$ws.Range("B26").Value = "20 Pin Female header socket"
$ws.Range("C26").Value = "20 pin socket"
$ws.Range("E26").Value = "20 Pin Female header sockets"
$ws.Range("F26").Value = "No"
$ws.Range("G26").Value = ""
$ws.Range("H26").Value = ""

# --- Row 36: mono Thonkiconn part number correction ---
$ws.Range("C36").Value = " PJ301M-12 / PJ398SM"
$ws.Range("E36").Value = "mono Thonkiconn"

# --- Restore the user's last selection state recorded in the workbook ---
$ws.Range("E37").Select()
